$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3322939521381819
$ws.Range("C2").Value = 0.9936074032634841
$ws.Range("D2").Value = 0.4553638524317597
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 AdaBoostRegressor(learning_rate=0.5, n_estimators=100))])"
$ws.Range("G2").Value = 0.1237476138499915
$ws.Range("H2").Value = 0.992
